$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.569.37"
$ws.Range("E2").Value = "  -3.75%  "

$ws.Range("D3").Value = "2.399.08"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.52"
$ws.Range("E5").Value = "  -5.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.20"
$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  -2.45%  "

$ws.Range("D9").Value = "2.397.44"
$ws.Range("E9").Value = "  -4.10%  "

$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("E12").Value = "  -1.52%  "

$ws.Range("E13").Value = "  -10.21%  "

$ws.Range("D14").Value = "2.798.76"
$ws.Range("E14").Value = "  -4.58%  "

$ws.Range("D15").Value = "56.435.06"
$ws.Range("E15").Value = "  -3.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.67"
$ws.Range("E16").Value = "  -2.52%  "

$ws.Range("E17").Value = "  -2.79%  "

$ws.Range("D18").Value = "2.390.37"
$ws.Range("E18").Value = "  -4.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.25"
$ws.Range("E19").Value = "  -3.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.89"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("E21").Value = "  -4.15%  "

$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.63"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("D26").Value = "2.495.61"
$ws.Range("E26").Value = "  -4.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.379"
$ws.Range("E27").Value = "  -6.70%  "

$ws.Range("E28").Value = "  -4.85%  "

$ws.Range("E29").Value = "  -2.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.52"
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("E32").Value = "  -5.14%  "

$ws.Range("E33").Value = "  -1.64%  "

$ws.Range("E34").Value = "  -6.25%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.82"
$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.77"
$ws.Range("E39").Value = "  -4.39%  "

$ws.Range("E40").Value = "  -1.26%  "

$ws.Range("E41").Value = "  -4.40%  "

$ws.Range("E42").Value = "  -4.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "132.51"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("E44").Value = "  -2.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.89"
$ws.Range("E45").Value = "  -3.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "255.55"
$ws.Range("E46").Value = "  -6.62%  "

$ws.Range("E47").Value = "  -3.34%  "

$ws.Range("E48").Value = "  -3.26%  "

$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("E50").Value = "  -4.09%  "

$ws.Range("E51").Value = "  -4.23%  "
